# Append: 2025-10-16 18:32 JST
# Update the "取得日時" timestamp for all data rows, and swap the
# title/URL of rows 6 and 7 (a re-ordering / de-dup fix upstream).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-16 18:32:53"

# Update column A (取得日時) for data rows 2..15
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $newTimestamp
}

# Swap title (B) and URL (F) between row 6 and row 7
$titleRow6 = $ws.Cells.Item(6, 2).Value2
$titleRow7 = $ws.Cells.Item(7, 2).Value2
$urlRow6 = $ws.Cells.Item(6, 6).Value2
$urlRow7 = $ws.Cells.Item(7, 6).Value2

$ws.Cells.Item(6, 2).Value2 = $titleRow7
$ws.Cells.Item(7, 2).Value2 = $titleRow6

$ws.Cells.Item(6, 6).Value2 = $urlRow7
$ws.Cells.Item(7, 6).Value2 = $urlRow6
